# Update the "想去人数" (interested-count) values (column F) on the
# "展览" sheet and the "全部类型" sheet, matching the rows that changed
# in the source data refresh (commit "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3179
$ws1.Range("F5").Value  = 2212
$ws1.Range("F6").Value  = 331
$ws1.Range("F8").Value  = 1063
$ws1.Range("F9").Value  = 1022
$ws1.Range("F10").Value = 246
$ws1.Range("F11").Value = 466
$ws1.Range("F12").Value = 1161
$ws1.Range("F16").Value = 7866
$ws1.Range("F17").Value = 345
$ws1.Range("F18").Value = 2469
$ws1.Range("F20").Value = 233
$ws1.Range("F23").Value = 540
$ws1.Range("F27").Value = 1523
$ws1.Range("F28").Value = 9
$ws1.Range("F30").Value = 1663
$ws1.Range("F38").Value = 179
$ws1.Range("F41").Value = 220

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 3179
$ws4.Range("F7").Value  = 2212
$ws4.Range("F8").Value  = 331
$ws4.Range("F10").Value = 1063
$ws4.Range("F12").Value = 1022
$ws4.Range("F13").Value = 246
$ws4.Range("F14").Value = 466
$ws4.Range("F15").Value = 1161
$ws4.Range("F19").Value = 7866
$ws4.Range("F20").Value = 345
$ws4.Range("F21").Value = 2469
$ws4.Range("F24").Value = 233
$ws4.Range("F27").Value = 540
$ws4.Range("F31").Value = 1523
$ws4.Range("F32").Value = 9
$ws4.Range("F34").Value = 1663
$ws4.Range("F42").Value = 179
$ws4.Range("F48").Value = 220
